# Changes of 21st June 2022
#
# The FedEx rate-verification sheet previously showed several rows as
# FAIL because the recorded "ActualRate" (column E) did not match the
# "ExpectedRate" (column D). Re-running / re-validating the rates shows
# they now agree, so ActualRate is updated to mirror ExpectedRate and the
# Result (column F) flips to PASS for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows whose ActualRate now matches ExpectedRate (and so Result -> PASS).
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 30)

foreach ($r in $rows) {
    $expectedCell = $ws.Cells.Item($r, 4)
    $actualCell   = $ws.Cells.Item($r, 5)
    $resultCell   = $ws.Cells.Item($r, 6)

    # Grab the displayed text of ExpectedRate (handles both numeric cells
    # formatted as currency, e.g. 19.04 -> "$19.04", and cells that are
    # already stored as literal currency text, e.g. "$473.23").
    $expectedText = $expectedCell.Text.Trim()

    # Force the destination cell to plain text first so Excel doesn't
    # "smart" re-parse the leading "$" back into a formatted number -
    # we want it stored the same way the other ActualRate text cells are
    # (shared string, no number formatting applied to the cell itself).
    $actualCell.NumberFormat = "@"
    $actualCell.Value = $expectedText
    $actualCell.ClearFormats()

    $resultCell.Value = "PASS"
}
